# Update course Excel file: change "department" value from
# "FACULTY OF ENGLISH" to "English" for the course rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# Column C holds "department" (row 1 is the header row).
# Rows 2 and 3 currently contain "FACULTY OF ENGLISH" - update them to "English".
$ws.Range("C2").Value = "English"
$ws.Range("C3").Value = "English"

$wb.Save()
